$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the sire/dam (F/G) values for rows 2, 5, 6, 7, 10 (keep header row 1)
$ws.Range("F2:G2").ClearContents()
$ws.Range("F5:G5").ClearContents()
$ws.Range("F6:G6").ClearContents()
$ws.Range("F7:G7").ClearContents()
$ws.Range("F10:G10").ClearContents()

# Update the selected cell/range shown in the saved file
$ws.Range("H21").Select()
